$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.578.74'
$ws.Range('E2').Value = '  +3.65%  '
$ws.Range('D3').Value = '1.916.96'
$ws.Range('E3').Value = '  +1.96%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.696'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.47%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.06'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.53%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '58.82'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +10.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.366'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.29%  '
$ws.Range('E11').Value = '  +3.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0998'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.78%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.57'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +8.13%  '
$ws.Range('E14').Value = '  +4.36%  '
$ws.Range('D15').Value = '2.196.95'
$ws.Range('E15').Value = '  +1.97%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.13'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.84%  '
$ws.Range('D17').Value = '1.917.22'
$ws.Range('E17').Value = '  +1.53%  '
$ws.Range('D18').Value = '36.556.67'
$ws.Range('E18').Value = '  +3.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.17'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.92%  '
$ws.Range('D20').Value = '0.0₃0861'
$ws.Range('E20').Value = '  +5.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '251.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.23'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.68%  '
$ws.Range('E23').Value = '  +5.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.68'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.61%  '
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.19'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.55'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.80'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.76'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.67%  '
$ws.Range('E30').Value = '  +1.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.55'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0608'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.86%  '
$ws.Range('E33').Value = '  +0.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.34'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.38%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0847'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +17.28%  '
$ws.Range('E37').Value = '  -13.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.873'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.84'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +47.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.02'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.69%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '106.63'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +11.17%  '
$ws.Range('E42').Value = '  +5.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.16'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.10'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.95%  '
$ws.Range('D45').Value = '1.337.54'
$ws.Range('E45').Value = '  +2.58%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.36'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.52'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.80%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0814'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.33%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.80'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.42'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '42.93'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.21%  '
